# Checkpoint III presentation - apply commit "Put names and group numbers in presentation."
$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Slide-number field placeholder cache text: "<nº>" -> "<#>"
#    (handout master, notes master, and the "Title and Content" slide layout)
# ---------------------------------------------------------------------
$numChar = "‹#›"

$handout = $p.HandoutMaster
for ($i = 1; $i -le $handout.Shapes.Count; $i++) {
    $shp = $handout.Shapes.Item($i)
    if ($shp.HasTextFrame -and $shp.Name -like "Slide Number Placeholder*") {
        $shp.TextFrame.TextRange.Text = $numChar
    }
}

$notesMaster = $p.NotesMaster
for ($i = 1; $i -le $notesMaster.Shapes.Count; $i++) {
    $shp = $notesMaster.Shapes.Item($i)
    if ($shp.HasTextFrame -and $shp.Name -like "Slide Number Placeholder*") {
        $shp.TextFrame.TextRange.Text = $numChar
    }
}

# Slide 3 uses the "Title and Content" layout (slideLayout2.xml) that also
# carries a slide-number placeholder field.
$layout = $p.Slides.Item(3).CustomLayout
for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
    $shp = $layout.Shapes.Item($i)
    if ($shp.HasTextFrame -and $shp.Name -like "Slide Number Placeholder*") {
        $shp.TextFrame.TextRange.Text = $numChar
    }
}

# ---------------------------------------------------------------------
# 2) Slide 1: group code box "GX-A/T" -> "G13-A"
# ---------------------------------------------------------------------
$slide1 = $p.Slides.Item(1)
$codeShape = $slide1.Shapes.Item(2)
$codeShape.TextFrame.TextRange.Text = "G13-A"

# ---------------------------------------------------------------------
# 3) Slide 1: resize/reposition the names box, and replace the
#    "Number - Name" placeholder lines with the real group members.
# ---------------------------------------------------------------------
$namesShape = $slide1.Shapes.Item(3)
$namesShape.Top = 4365104 / 12700
$namesShape.Width = 2664296 / 12700
$namesShape.Height = 2492896 / 12700

$tr = $namesShape.TextFrame.TextRange

function Set-ParaFullText($paragraph, [string]$newText) {
    # Each paragraph currently reads "Number - Name" as three runs:
    #   "Number" (err=1) + " - " (clean) + "Name" (err=1)
    # Drop the leading "Number" run and the trailing "Name" run, keeping
    # the clean middle run's formatting, then put the full text into it.
    $numberRun = $paragraph.Characters(1, 6)
    $numberRun.Text = ""

    $len2 = $paragraph.Characters().Count
    $nameRun = $paragraph.Characters($len2 - 3, 4)
    $nameRun.Text = ""

    $remaining = $paragraph.Characters()
    $remaining.Text = $newText
}

Set-ParaFullText $tr.Paragraphs(1) "83463 – Francisco Campaniço"
Set-ParaFullText $tr.Paragraphs(2) "83482 – João Rafael"
Set-ParaFullText $tr.Paragraphs(3) "83558 – Rodrigo Oliveira"
